$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.339.71'
$ws.Range('E2').Value = '  +0.79%  '
$ws.Range('D3').Value = '3.787.74'
$ws.Range('E3').Value = '  +1.12%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '603.85'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '164.96'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -2.24%  '
$ws.Range('D7').Value = '3.780.24'
$ws.Range('E7').Value = '  +1.02%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  +0.82%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.33'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.39%  '
$ws.Range('E12').Value = '  -0.34%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '37.49'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -2.17%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000247'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.20%  '
$ws.Range('D15').Value = '4.417.88'
$ws.Range('E15').Value = '  +0.95%  '
$ws.Range('D16').Value = '3.790.07'
$ws.Range('E16').Value = '  +1.28%  '
$ws.Range('D17').Value = '69.368.20'
$ws.Range('E17').Value = '  +0.82%  '
$ws.Range('E18').Value = '  +1.92%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.60'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +3.06%  '
$ws.Range('E20').Value = '  -1.10%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.27'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +4.37%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '493.04'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.71%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.724'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.59%  '
$ws.Range('E24').Value = '  -1.92%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.84'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.67%  '
$ws.Range('E26').Value = '  -2.13%  '
$ws.Range('E27').Value = '  -0.17%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.14'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.57%  '
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.99'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.12%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.14'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +2.87%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.42'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -3.95%  '
$ws.Range('B33').Value = 'WrappedeETH'
$ws.Range('C33').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D33').Value = '3.932.13'
$ws.Range('E33').Value = '  +1.03%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '31.87'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.17%  '
$ws.Range('D35').Value = '3.731.34'
$ws.Range('E35').Value = '  +1.34%  '
$ws.Range('E36').Value = '  -0.60%  '
$ws.Range('B37').Value = 'Mantle'
$ws.Range('C37').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.02'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.03%  '
$ws.Range('B38').Value = 'Filecoin'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.95'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.56%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.139'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +5.01%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.999'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.09%  '
$ws.Range('E41').Value = '  +0.17%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.09'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +4.79%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '48.73'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.20%  '
$ws.Range('E44').Value = '  +1.20%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '423.33'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -3.57%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.43'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.46%  '
$ws.Range('E47').Value = '  +0.01%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '40.04'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.18%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '141.87'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.03%  '
$ws.Range('B50').Value = 'ONDO'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.31'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +7.64%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = '2.812.45'
$ws.Range('E51').Value = '  +0.83%  '
